$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.94338166666667
$ws.Range("H2").Value = 47.830145
$ws.Range("I2").Value = 0.09439898960585606
$ws.Range("J2").Value = 0.09439898960585606
$ws.Range("M2").Value = 57.353185
$ws.Range("N2").Value = 172.059555
$ws.Range("O2").Value = 0.2415415519323779
$ws.Range("P2").Value = 0.2415415519323779
$ws.Range("Q2").Value = 914.4037182539416
$ws.Range("R2").Value = 8229.633464285474
$ws.Range("S2").Value = 0.02280127845024689
$ws.Range("T2").Value = 0.02280127845024689
$ws.Range("G3").Value = 15.94338166666667
$ws.Range("H3").Value = 47.830145
$ws.Range("I3").Value = 0.09439898960585606
$ws.Range("J3").Value = 0.09439898960585606
$ws.Range("O3").Value = 0.1738483040615689
$ws.Range("P3").Value = 0.1738483040615689
$ws.Range("Q3").Value = 658.137427594839
$ws.Range("R3").Value = 5923.236848353551
$ws.Range("S3").Value = 0.01641110424810375
$ws.Range("T3").Value = 0.01641110424810375
$ws.Range("G4").Value = 15.94338166666667
$ws.Range("H4").Value = 47.830145
$ws.Range("I4").Value = 0.09439898960585606
$ws.Range("J4").Value = 0.09439898960585606
$ws.Range("M4").Value = 17.547551
$ws.Range("N4").Value = 52.642653
$ws.Range("O4").Value = 0.07390108676880894
$ws.Range("P4").Value = 0.07390108676880894
$ws.Range("Q4").Value = 279.7673029082984
$ws.Range("R4").Value = 2517.905726174685
$ws.Range("S4").Value = 0.006976187921750262
$ws.Range("T4").Value = 0.006976187921750262
$ws.Range("G5").Value = 15.94338166666667
$ws.Range("H5").Value = 47.830145
$ws.Range("I5").Value = 0.09439898960585606
$ws.Range("J5").Value = 0.09439898960585606
$ws.Range("M5").Value = 121.2660546666667
$ws.Range("N5").Value = 363.798164
$ws.Range("O5").Value = 0.5107090572372441
$ws.Range("P5").Value = 0.5107090572372442
$ws.Range("Q5").Value = 1933.390992761531
$ws.Range("R5").Value = 17400.51893485378
$ws.Range("S5").Value = 0.04821041898575516
$ws.Range("T5").Value = 0.04821041898575517
$ws.Range("I6").Value = 0.07777357403446172
$ws.Range("J6").Value = 0.07777357403446171
$ws.Range("M6").Value = 57.353185
$ws.Range("N6").Value = 172.059555
$ws.Range("O6").Value = 0.2415415519323779
$ws.Range("P6").Value = 0.2415415519323779
$ws.Range("Q6").Value = 753.36023802736
$ws.Range("R6").Value = 6780.242142246239
$ws.Range("S6").Value = 0.01878554977161158
$ws.Range("T6").Value = 0.01878554977161158
$ws.Range("I7").Value = 0.07777357403446172
$ws.Range("J7").Value = 0.07777357403446171
$ws.Range("O7").Value = 0.1738483040615689
$ws.Range("P7").Value = 0.1738483040615689
$ws.Range("S7").Value = 0.01352080394669804
$ws.Range("T7").Value = 0.01352080394669804
$ws.Range("I8").Value = 0.07777357403446172
$ws.Range("J8").Value = 0.07777357403446171
$ws.Range("M8").Value = 17.547551
$ws.Range("N8").Value = 52.642653
$ws.Range("O8").Value = 0.07390108676880894
$ws.Range("P8").Value = 0.07390108676880894
$ws.Range("Q8").Value = 230.495084068256
$ws.Range("R8").Value = 2074.455756614304
$ws.Range("S8").Value = 0.005747551643041141
$ws.Range("T8").Value = 0.00574755164304114
$ws.Range("I9").Value = 0.07777357403446172
$ws.Range("J9").Value = 0.07777357403446171
$ws.Range("M9").Value = 121.2660546666667
$ws.Range("N9").Value = 363.798164
$ws.Range("O9").Value = 0.5107090572372441
$ws.Range("P9").Value = 0.5107090572372442
$ws.Range("Q9").Value = 1592.884925367595
$ws.Range("R9").Value = 14335.96432830835
$ws.Range("S9").Value = 0.03971966867311096
$ws.Range("T9").Value = 0.03971966867311096
$ws.Range("G10").Value = 1.134259333333333
$ws.Range("H10").Value = 3.402778
$ws.Range("I10").Value = 0.006715823358951466
$ws.Range("J10").Value = 0.006715823358951466
$ws.Range("M10").Value = 57.353185
$ws.Range("N10").Value = 172.059555
$ws.Range("O10").Value = 0.2415415519323779
$ws.Range("P10").Value = 0.2415415519323779
$ws.Range("Q10").Value = 65.05338538264333
$ws.Range("R10").Value = 585.48046844379
$ws.Range("S10").Value = 0.001622150396624852
$ws.Range("T10").Value = 0.001622150396624852
$ws.Range("G11").Value = 1.134259333333333
$ws.Range("H11").Value = 3.402778
$ws.Range("I11").Value = 0.006715823358951466
$ws.Range("J11").Value = 0.006715823358951466
$ws.Range("O11").Value = 0.1738483040615689
$ws.Range("P11").Value = 0.1738483040615689
$ws.Range("Q11").Value = 46.82184341269112
$ws.Range("R11").Value = 421.3965907142201
$ws.Range("S11").Value = 0.001167534501330782
$ws.Range("T11").Value = 0.001167534501330781
$ws.Range("G12").Value = 1.134259333333333
$ws.Range("H12").Value = 3.402778
$ws.Range("I12").Value = 0.006715823358951466
$ws.Range("J12").Value = 0.006715823358951466
$ws.Range("M12").Value = 17.547551
$ws.Range("N12").Value = 52.642653
$ws.Range("O12").Value = 0.07390108676880894
$ws.Range("P12").Value = 0.07390108676880894
$ws.Range("Q12").Value = 19.90347349889267
$ws.Range("R12").Value = 179.131261490034
$ws.Range("S12").Value = 0.0004963066447738662
$ws.Range("T12").Value = 0.0004963066447738662
$ws.Range("G13").Value = 1.134259333333333
$ws.Range("H13").Value = 3.402778
$ws.Range("I13").Value = 0.006715823358951466
$ws.Range("J13").Value = 0.006715823358951466
$ws.Range("M13").Value = 121.2660546666667
$ws.Range("N13").Value = 363.798164
$ws.Range("O13").Value = 0.5107090572372441
$ws.Range("P13").Value = 0.5107090572372442
$ws.Range("Q13").Value = 137.5471543221769
$ws.Range("R13").Value = 1237.924388899592
$ws.Range("S13").Value = 0.003429831816221965
$ws.Range("T13").Value = 0.003429831816221966
$ws.Range("G14").Value = 138.6804656666667
$ws.Range("H14").Value = 416.041397
$ws.Range("I14").Value = 0.8211116130007308
$ws.Range("J14").Value = 0.8211116130007308
$ws.Range("M14").Value = 57.353185
$ws.Range("N14").Value = 172.059555
$ws.Range("O14").Value = 0.2415415519323779
$ws.Range("P14").Value = 0.2415415519323779
$ws.Range("Q14").Value = 7953.766403266481
$ws.Range("R14").Value = 71583.89762939833
$ws.Range("S14").Value = 0.1983325733138946
$ws.Range("T14").Value = 0.1983325733138946
$ws.Range("G15").Value = 138.6804656666667
$ws.Range("H15").Value = 416.041397
$ws.Range("I15").Value = 0.8211116130007308
$ws.Range("J15").Value = 0.8211116130007308
$ws.Range("O15").Value = 0.1738483040615689
$ws.Range("P15").Value = 0.1738483040615689
$ws.Range("Q15").Value = 5724.682933629892
$ws.Range("R15").Value = 51522.14640266904
$ws.Range("S15").Value = 0.1427488613654364
$ws.Range("T15").Value = 0.1427488613654363
$ws.Range("G16").Value = 138.6804656666667
$ws.Range("H16").Value = 416.041397
$ws.Range("I16").Value = 0.8211116130007308
$ws.Range("J16").Value = 0.8211116130007308
$ws.Range("M16").Value = 17.547551
$ws.Range("N16").Value = 52.642653
$ws.Range("O16").Value = 0.07390108676880894
$ws.Range("P16").Value = 0.07390108676880894
$ws.Range("Q16").Value = 2433.502543989583
$ws.Range("R16").Value = 21901.52289590624
$ws.Range("S16").Value = 0.06068104055924367
$ws.Range("T16").Value = 0.06068104055924367
$ws.Range("G17").Value = 138.6804656666667
$ws.Range("H17").Value = 416.041397
$ws.Range("I17").Value = 0.8211116130007308
$ws.Range("J17").Value = 0.8211116130007308
$ws.Range("M17").Value = 121.2660546666667
$ws.Range("N17").Value = 363.798164
$ws.Range("O17").Value = 0.5107090572372441
$ws.Range("P17").Value = 0.5107090572372442
$ws.Range("Q17").Value = 16817.23293073279
$ws.Range("R17").Value = 151355.0963765951
$ws.Range("S17").Value = 0.419349137762156
$ws.Range("T17").Value = 0.4193491377621562
